$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (129) down into the
# two new rows (130 and 131) so the new cells inherit the same styles
# (bold/border index style on column A, datetime numFmt on column E).
$ws.Range("A129:V129").Copy()
$ws.Range("A130:V131").PasteSpecial(-4122)

# ---- Row 130 (new match #129) ----
$ws.Cells.Item(130, 1).Value  = 129
$ws.Cells.Item(130, 2).Value  = "poland"
$ws.Cells.Item(130, 3).Value  = "ekstraklasa"
$ws.Cells.Item(130, 4).Value  = "2023-2024"
$ws.Cells.Item(130, 5).Value  = 45242.52083333334
$ws.Cells.Item(130, 6).Value  = "Stal Mielec"
$ws.Cells.Item(130, 7).Value  = 2
$ws.Cells.Item(130, 8).Value  = "Gornik Zabrze"
$ws.Cells.Item(130, 9).Value  = 1
$ws.Cells.Item(130, 10).Value = 2.69
$ws.Cells.Item(130, 11).Value = "06/11/2023 19:12"
$ws.Cells.Item(130, 12).Value = 3.13
$ws.Cells.Item(130, 13).Value = "12/11/2023 12:20"
$ws.Cells.Item(130, 14).Value = 3.19
$ws.Cells.Item(130, 15).Value = "06/11/2023 19:12"
$ws.Cells.Item(130, 16).Value = 3.33
$ws.Cells.Item(130, 17).Value = "12/11/2023 12:20"
$ws.Cells.Item(130, 18).Value = 2.67
$ws.Cells.Item(130, 19).Value = "06/11/2023 19:12"
$ws.Cells.Item(130, 20).Value = 2.42
$ws.Cells.Item(130, 21).Value = "12/11/2023 12:20"
$ws.Cells.Item(130, 22).Value = "https://www.betexplorer.com/football/poland/ekstraklasa/stal-mielec-gornik-zabrze/WdKYTZ7j/"

# ---- Row 131 (new match #130) ----
$ws.Cells.Item(131, 1).Value  = 130
$ws.Cells.Item(131, 2).Value  = "poland"
$ws.Cells.Item(131, 3).Value  = "ekstraklasa"
$ws.Cells.Item(131, 4).Value  = "2023-2024"
$ws.Cells.Item(131, 5).Value  = 45242.625
$ws.Cells.Item(131, 6).Value  = "Pogon Szczecin"
$ws.Cells.Item(131, 7).Value  = 1
$ws.Cells.Item(131, 8).Value  = "Rakow"
$ws.Cells.Item(131, 9).Value  = 1
$ws.Cells.Item(131, 10).Value = 2.21
$ws.Cells.Item(131, 11).Value = "05/11/2023 15:13"
$ws.Cells.Item(131, 12).Value = 2.61
$ws.Cells.Item(131, 13).Value = "12/11/2023 14:53"
$ws.Cells.Item(131, 14).Value = 3.61
$ws.Cells.Item(131, 15).Value = "05/11/2023 15:13"
$ws.Cells.Item(131, 16).Value = 3.41
$ws.Cells.Item(131, 17).Value = "12/11/2023 14:50"
$ws.Cells.Item(131, 18).Value = 3.03
$ws.Cells.Item(131, 19).Value = "05/11/2023 15:13"
$ws.Cells.Item(131, 20).Value = 2.81
$ws.Cells.Item(131, 21).Value = "12/11/2023 14:53"
$ws.Cells.Item(131, 22).Value = "https://www.betexplorer.com/football/poland/ekstraklasa/pogon-szczecin-rakow-czestochowa/zRY2ZepT/"
